$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the pre-existing hyperlink on E2 (its URL/hyperlink is moving to column G)
$ws.Hyperlinks.Item(1).Delete()

# New header cells for Artist / Publication Date, and URL header shifts out to column G
$ws.Range("E1").Value = "Artist"
$ws.Range("F1").Value = "Publication Date"
$ws.Range("G1").Value = "URL"

# Row 2 (Sample Tomato Sound): Artist added in E, URL + hyperlink moved to G
$ws.Range("E2").Value = "Michael"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "www.FakeURL.com/sampleTomato"

# Row 3 (50 Cal Casing): Artist + Publication Date added, URL + hyperlink moved to G
$ws.Range("E3").Value = "Ghost Rider"
$ws.Range("F3").Value = 40005
$ws.Range("F3").NumberFormat = "mm-dd-yy"
$ws.Range("G3").Value = "http://soundbible.com/1927-50-Cal-CasingX3.html"

# Row 4 (new Footsteps On Cement sound)
$ws.Range("D4").Value = ".mp3"
$ws.Range("E4").Value = "Tim Fryer"
$ws.Range("F4").Value = 41346
$ws.Range("F4").NumberFormat = "mm-dd-yy"
$ws.Range("G4").Value = "http://soundbible.com/2057-Footsteps-On-Cement.html"

# Hyperlinks for the URL column
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.fakeurl.com/sampleTomato", "", "", "www.FakeURL.com/sampleTomato")
$ws.Range("G2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G3"), "http://soundbible.com/1927-50-Cal-CasingX3.html")
$ws.Range("G3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("G4"), "http://soundbible.com/2057-Footsteps-On-Cement.html")
$ws.Range("G4").Style = "Hyperlink"

# Column widths for new columns F and G
$ws.Columns.Item(6).ColumnWidth = 18.08984375
$ws.Columns.Item(7).ColumnWidth = 53.453125

# View state: scroll so column E is leftmost, select G5 (first blank row below the data)
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("G5").Select()
